$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute('2024-06-08 Saturday', $true, $true, $false, $false, $false, $true, 1, $false, '2024-06-09 Sunday', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('682÷3=', $true, $true, $false, $false, $false, $true, 1, $false, '897÷8=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('357÷9=', $true, $true, $false, $false, $false, $true, 1, $false, '610÷3=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('115÷5=', $true, $true, $false, $false, $false, $true, 1, $false, '337÷8=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('773÷6=', $true, $true, $false, $false, $false, $true, 1, $false, '703÷3=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('853÷8=', $true, $true, $false, $false, $false, $true, 1, $false, '822÷5=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('276÷2=', $true, $true, $false, $false, $false, $true, 1, $false, '512÷5=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('279÷6=', $true, $true, $false, $false, $false, $true, 1, $false, '651÷7=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('556÷4=', $true, $true, $false, $false, $false, $true, 1, $false, '195÷5=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('289÷8=', $true, $true, $false, $false, $false, $true, 1, $false, '711÷8=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('210÷8=', $true, $true, $false, $false, $false, $true, 1, $false, '994÷9=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('452÷3=', $true, $true, $false, $false, $false, $true, 1, $false, '825÷2=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('200÷4=', $true, $true, $false, $false, $false, $true, 1, $false, '848÷2=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('401÷4=', $true, $true, $false, $false, $false, $true, 1, $false, '163÷5=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('884÷9=', $true, $true, $false, $false, $false, $true, 1, $false, '121÷6=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('776÷5=', $true, $true, $false, $false, $false, $true, 1, $false, '558÷3=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('923÷6=', $true, $true, $false, $false, $false, $true, 1, $false, '287÷4=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('854÷2=', $true, $true, $false, $false, $false, $true, 1, $false, '226÷8=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('853÷6=', $true, $true, $false, $false, $false, $true, 1, $false, '377÷6=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('174÷3=', $true, $true, $false, $false, $false, $true, 1, $false, '519÷8=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('453÷5=', $true, $true, $false, $false, $false, $true, 1, $false, '391÷7=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('236÷6=', $true, $true, $false, $false, $false, $true, 1, $false, '637÷2=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('684÷4=', $true, $true, $false, $false, $false, $true, 1, $false, '844÷3=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('578÷8=', $true, $true, $false, $false, $false, $true, 1, $false, '680÷2=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('938÷9=', $true, $true, $false, $false, $false, $true, 1, $false, '238÷8=', 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute('721÷2=', $true, $true, $false, $false, $false, $true, 1, $false, '366÷8=', 2) | Out-Null
